# Renamed U15 to U16.
# Update the command names in the "XMOS->STM32" sheet that referred to the
# U15 IO expander so that they now refer to U16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS->STM32")

$ws.Range("A42").Value = "Set U16 outputs"
$ws.Range("A41").Value = "Get U16 outputs"
$ws.Range("A40").Value = "Set U16 output pin X"
$ws.Range("A39").Value = "Get U16 output pin X"

$ws.Range("A38").Select()
